# Updating weight measurements: add a new weekly measurement row and
# switch the "Weight (Stone)" column to a text number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new measurement row (row 5) -----------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 43115
$ws.Range("C5").Value = 15.13
$ws.Range("D5").Value = 101.1
$ws.Range("E5").Value = 222.8
$ws.Range("F5").Formula = "=E5-E4"
$ws.Range("G5").Formula = "=E5-210"

# --- Re-format column C (Weight (Stone)) as text, header through new row
$ws.Range("C1:C5").NumberFormat = "@"

# --- Move the selection like the original author left it -------------
$ws.Range("B6").Select() | Out-Null
